$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Regression coefficients for the "Crisis and Credit Allocation" table.
# Values that are plain numerics (no trailing significance stars) need a
# leading apostrophe so Excel keeps storing them as text, matching the
# original workbook's convention of keeping every coefficient as a text
# cell (even when it parses as a plain number) instead of as a number.

# A Lag row
$ws.Range("B2").Value = "'0.25"
$ws.Range("C2").Value = "0.4***"
$ws.Range("D2").Value = "-11.46***"

# C Lag row
$ws.Range("B3").Value = "'0.12"
$ws.Range("C3").Value = "-0.57***"
$ws.Range("D3").Value = "'10.16"

# LF Lag row
$ws.Range("B4").Value = "-0.1*"
$ws.Range("C4").Value = "'0.02"
$ws.Range("D4").Value = "1.48*"
